$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "21.760.01"
$ws.Range("E2").Value = "  +5.87%  "
$ws.Range("D3").Value = "1.576.00"
$ws.Range("E3").Value = "  +6.67%  "
$ws.Range("D4").Value = "'0.9944"
$ws.Range("E4").Value = "  -1.24%  "
$ws.Range("D5").Value = "'0.9692"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").Value = "'284.49"
$ws.Range("E6").Value = "  +2.74%  "
$ws.Range("D7").Value = "'0.3676"
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("D8").Value = "'0.3261"
$ws.Range("E8").Value = "  +6.55%  "
$ws.Range("D9").Value = "'1.133"
$ws.Range("E9").Value = "  +6.90%  "
$ws.Range("D10").Value = "'41.06"
$ws.Range("E10").Value = "  +3.29%  "
$ws.Range("D11").Value = "'0.07040"
$ws.Range("E11").Value = "  +6.17%  "
$ws.Range("D12").Value = "'0.9884"
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("D13").Value = "'20.08"
$ws.Range("E13").Value = "  +10.15%  "
$ws.Range("D14").Value = "'5.784"
$ws.Range("E14").Value = "  +5.56%  "
$ws.Range("D15").Value = "'6.475"
$ws.Range("E15").Value = "  +4.81%  "
$ws.Range("D16").Value = "'0.00001067"
$ws.Range("E16").Value = "  +3.52%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "1.561.63"
$ws.Range("E17").Value = "  +5.80%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'0.9677"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").Value = "'0.06158"
$ws.Range("E19").Value = "  +4.21%  "
$ws.Range("D20").Value = "'73.72"
$ws.Range("E20").Value = "  +6.21%  "
$ws.Range("D21").Value = "'15.92"
$ws.Range("E21").Value = "  +9.81%  "
$ws.Range("D22").Value = "'5.811"
$ws.Range("E22").Value = "  +6.10%  "
$ws.Range("E23").Value = "  +5.07%  "
$ws.Range("B24").Value = "WrappedBTC"
$ws.Range("C24").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D24").Value = "21.669.93"
$ws.Range("E24").Value = "  +5.26%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'2.320"
$ws.Range("E25").Value = "  +3.30%  "
$ws.Range("D26").Value = "'2.416"
$ws.Range("E26").Value = "  +12.91%  "
$ws.Range("D27").Value = "'148.35"
$ws.Range("E27").Value = "  +4.77%  "
$ws.Range("E28").Value = "  +5.26%  "
$ws.Range("D29").Value = "1.728.87"
$ws.Range("E29").Value = "  +5.94%  "
$ws.Range("D30").Value = "'119.55"
$ws.Range("E30").Value = "  +5.04%  "
$ws.Range("D31").Value = "'4.051"
$ws.Range("E31").Value = "  +3.50%  "
$ws.Range("D32").Value = "'0.8989"
$ws.Range("E32").Value = "  +9.34%  "
$ws.Range("D33").Value = "'5.341"
$ws.Range("E33").Value = "  +7.08%  "
$ws.Range("D34").Value = "'0.08135"
$ws.Range("E34").Value = "  +2.34%  "
$ws.Range("D35").Value = "'1.553"
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("D36").Value = "'5.066"
$ws.Range("E36").Value = "  +6.87%  "
$ws.Range("D37").Value = "'1.234"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "'11.45"
$ws.Range("E38").Value = "  +9.26%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.05989"
$ws.Range("E39").Value = "  +3.92%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'8.174"
$ws.Range("E40").Value = "  +7.29%  "
$ws.Range("D41").Value = "'0.02161"
$ws.Range("E41").Value = "  +5.91%  "
$ws.Range("D42").Value = "'0.2007"
$ws.Range("E42").Value = "  +6.59%  "
$ws.Range("D43").Value = "'0.9672"
$ws.Range("E43").Value = "  +0.93%  "
$ws.Range("D44").Value = "'0.5754"
$ws.Range("E44").Value = "  +8.77%  "
$ws.Range("D45").Value = "'12.73"
$ws.Range("E45").Value = "  +5.02%  "
$ws.Range("D46").Value = "'3.595"
$ws.Range("E46").Value = "  +2.52%  "
$ws.Range("D47").Value = "'0.5639"
$ws.Range("E47").Value = "  +8.43%  "
$ws.Range("D48").Value = "'123.77"
$ws.Range("E48").Value = "  +5.46%  "
$ws.Range("D49").Value = "'1.909"
$ws.Range("E49").Value = "  +7.19%  "
$ws.Range("D50").Value = "'0.06716"
$ws.Range("E50").Value = "  +3.88%  "
$ws.Range("D51").Value = "'71.20"
$ws.Range("E51").Value = "  +5.74%  "